# The framework's LoginData test data had a duplicate row (row 5) that
# mirrored row 2/4 (training@jalaacademy.com / jobprogram). Clear that
# leftover row's contents now that the framework part is complete,
# while leaving the cell formatting/styles in place.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Range("A5:B5").ClearContents()
